# Auto-generated edit script: updates computed profit/price columns
# (H, I, J, K, L, M, N) on several rows across the job sheets, per the
# commit's scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 183511.64
$ws.Range("J17").Value = 183511.64
$ws.Range("L17").Value = 550534.92
$ws.Range("N17").Value = -550870.92

$ws.Range("H51").Value = 11888.581
$ws.Range("J51").Value = 11777.931
$ws.Range("L51").Value = 11777.931
$ws.Range("N51").Value = -12745.931

$ws.Range("H100").Value = 1393.7
$ws.Range("I100").Value = 1167.6923
$ws.Range("J100").Value = 1813.4286
$ws.Range("K100").Value = 1167.6923
$ws.Range("L100").Value = 1813.4286
$ws.Range("M100").Value = -626.6922999999999
$ws.Range("N100").Value = -2895.4286

$ws.Range("H103").Value = 1196.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1196.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 3589.5
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -4761.5

$ws.Range("H106").Value = 3557.4167
$ws.Range("I106").Value = 2749.5715
$ws.Range("J106").Value = 4688.4
$ws.Range("K106").Value = 2749.5715
$ws.Range("L106").Value = 4688.4
$ws.Range("M106").Value = -2118.5715
$ws.Range("N106").Value = -5950.4

$ws.Range("H112").Value = 41408.52
$ws.Range("I112").Value = 73029.57000000001
$ws.Range("J112").Value = 29758.658
$ws.Range("K112").Value = 219088.71
$ws.Range("L112").Value = 89275.974
$ws.Range("M112").Value = -217980.71
$ws.Range("N112").Value = -91491.974


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 52633930
$ws.Range("I61").Value = 58825616
$ws.Range("K61").Value = 58825616
$ws.Range("M61").Value = -58825404

$ws.Range("H102").Value = 999
$ws.Range("I102").Value = 999
$ws.Range("K102").Value = 999
$ws.Range("M102").Value = 623

$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws.Range("H122").Value = 1298.1052
$ws.Range("I122").Value = 1005.4545
$ws.Range("J122").Value = 3229.6
$ws.Range("K122").Value = 3016.3635
$ws.Range("L122").Value = 9688.799999999999
$ws.Range("M122").Value = -566.3635000000004
$ws.Range("N122").Value = -14588.8

$ws.Range("H136").Value = 52633930
$ws.Range("I136").Value = 58825616
$ws.Range("K136").Value = 176476848
$ws.Range("M136").Value = -176474298


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 33598.5
$ws.Range("I6").Value = 23647
$ws.Range("J6").Value = 39569.4
$ws.Range("K6").Value = 23647
$ws.Range("L6").Value = 39569.4
$ws.Range("M6").Value = -23534
$ws.Range("N6").Value = -39795.4

$ws.Range("H36").Value = 4875.125
$ws.Range("I36").Value = 1700
$ws.Range("J36").Value = 5933.5
$ws.Range("K36").Value = 1700
$ws.Range("L36").Value = 5933.5
$ws.Range("M36").Value = -1166
$ws.Range("N36").Value = -7001.5

$ws.Range("H82").Value = 42563.5
$ws.Range("I82").Value = 42563.5
$ws.Range("K82").Value = 42563.5
$ws.Range("M82").Value = -42180.5

$ws.Range("H85").Value = 42563.5
$ws.Range("I85").Value = 42563.5
$ws.Range("K85").Value = 42563.5
$ws.Range("M85").Value = -41237.5

$ws.Range("H94").Value = 205.83333

$ws.Range("H97").Value = 50284.8
$ws.Range("J97").Value = 89999.5
$ws.Range("L97").Value = 89999.5
$ws.Range("N97").Value = -91981.5

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 1833.6666
$ws.Range("J99").Value = 2499
$ws.Range("K99").Value = 1833.6666
$ws.Range("L99").Value = 2499
$ws.Range("M99").Value = -335.6666
$ws.Range("N99").Value = -5495

$ws.Range("H134").Value = 39617652
$ws.Range("I134").Value = 42918956
$ws.Range("K134").Value = 128756868
$ws.Range("M134").Value = -128754333


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 412.5
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 383.33334
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 383.33334
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -609.33334

$ws.Range("H31").Value = 14278.25
$ws.Range("I31").Value = 9850.833000000001
$ws.Range("K31").Value = 9850.833000000001
$ws.Range("M31").Value = -9555.833000000001

$ws.Range("H34").Value = 14278.25
$ws.Range("I34").Value = 9850.833000000001
$ws.Range("K34").Value = 9850.833000000001
$ws.Range("M34").Value = -9648.833000000001

$ws.Range("H62").Value = 3070.625
$ws.Range("I62").Value = 3115
$ws.Range("J62").Value = 2937.5
$ws.Range("K62").Value = 3115
$ws.Range("L62").Value = 2937.5
$ws.Range("M62").Value = -2491
$ws.Range("N62").Value = -4185.5

$ws.Range("H65").Value = 3070.625
$ws.Range("I65").Value = 3115
$ws.Range("J65").Value = 2937.5
$ws.Range("K65").Value = 15575
$ws.Range("L65").Value = 14687.5
$ws.Range("M65").Value = -12455
$ws.Range("N65").Value = -20927.5

$ws.Range("H134").Value = 8931554
$ws.Range("I134").Value = 9618327
$ws.Range("K134").Value = 28854981
$ws.Range("M134").Value = -28852446


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 50411.75
$ws.Range("I5").Value = 52933.42
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 158800.26
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -158688.26
$ws.Range("N5").Value = -7724

$ws.Range("H12").Value = 174.81818
$ws.Range("I12").Value = 246.71428
$ws.Range("K12").Value = 740.14284
$ws.Range("M12").Value = -567.14284

$ws.Range("H36").Value = 1002
$ws.Range("I36").Value = 1002
$ws.Range("K36").Value = 3006
$ws.Range("M36").Value = -2837

$ws.Range("H121").Value = 110689.73
$ws.Range("I121").Value = 250447.25
$ws.Range("J121").Value = 30828.285
$ws.Range("K121").Value = 751341.75
$ws.Range("L121").Value = 92484.855
$ws.Range("M121").Value = -750031.75
$ws.Range("N121").Value = -95104.855

$ws.Range("H135").Value = 50411.75
$ws.Range("I135").Value = 52933.42
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 476400.78
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -473865.78
$ws.Range("N135").Value = -27570


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2951.7917
$ws.Range("I80").Value = 2558.5833
$ws.Range("J80").Value = 3345
$ws.Range("K80").Value = 2558.5833
$ws.Range("L80").Value = 3345
$ws.Range("M80").Value = -1560.5833
$ws.Range("N80").Value = -5341

$ws.Range("H83").Value = 2951.7917
$ws.Range("I83").Value = 2558.5833
$ws.Range("J83").Value = 3345
$ws.Range("K83").Value = 12792.9165
$ws.Range("L83").Value = 16725
$ws.Range("M83").Value = -7800.916499999999
$ws.Range("N83").Value = -26709

$ws.Range("H87").Value = 69999
$ws.Range("J87").Value = 69999
$ws.Range("L87").Value = 69999
$ws.Range("N87").Value = -72495

$ws.Range("H90").Value = 69999
$ws.Range("J90").Value = 69999
$ws.Range("L90").Value = 209997
$ws.Range("N90").Value = -222477

$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -43744

$ws.Range("H97").Value = 644.8823
$ws.Range("J97").Value = 1099
$ws.Range("L97").Value = 1099
$ws.Range("N97").Value = -2091

$ws.Range("H107").Value = 1165.1666
$ws.Range("J107").Value = 1647.25
$ws.Range("L107").Value = 1647.25
$ws.Range("N107").Value = -5487.25

$ws.Range("H113").Value = 146984.72
$ws.Range("I113").Value = 251823.25
$ws.Range("K113").Value = 251823.25
$ws.Range("M113").Value = -249653.25

$ws.Range("H119").Value = 52000
$ws.Range("J119").Value = 52000
$ws.Range("L119").Value = 52000
$ws.Range("N119").Value = -61676

$ws.Range("H122").Value = 51644.48
$ws.Range("I122").Value = 65216.473
$ws.Range("K122").Value = 195649.419
$ws.Range("M122").Value = -193199.419

$ws.Range("H132").Value = 6252464
$ws.Range("I132").Value = 6581435
$ws.Range("J132").Value = 2014
$ws.Range("K132").Value = 19744305
$ws.Range("L132").Value = 6042
$ws.Range("M132").Value = -19741775
$ws.Range("N132").Value = -11102


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 211.71428
$ws.Range("I55").Value = 187
$ws.Range("K55").Value = 187
$ws.Range("M55").Value = -14

$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488

$ws.Range("H93").Value = 3549.5
$ws.Range("I93").Value = 3074.5
$ws.Range("K93").Value = 3074.5
$ws.Range("M93").Value = -1826.5

$ws.Range("H100").Value = 25003328
$ws.Range("I100").Value = 25003328
$ws.Range("K100").Value = 25003328
$ws.Range("M100").Value = -25002787

$ws.Range("H122").Value = 773968.1
$ws.Range("I122").Value = 773968.1
$ws.Range("K122").Value = 2321904.3
$ws.Range("M122").Value = -2319454.3

$ws.Range("H132").Value = 21825104
$ws.Range("I132").Value = 21825104
$ws.Range("K132").Value = 65475312
$ws.Range("M132").Value = -65472782


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3272.5454
$ws.Range("J96").Value = 2332.6667
$ws.Range("L96").Value = 2332.6667
$ws.Range("N96").Value = -5078.6667

$ws.Range("H122").Value = 3335600
$ws.Range("I122").Value = 3335600
$ws.Range("K122").Value = 10006800
$ws.Range("M122").Value = -10004350

$ws.Range("H132").Value = 21746112
$ws.Range("I132").Value = 31252776
$ws.Range("J132").Value = 16596.285
$ws.Range("K132").Value = 93758328
$ws.Range("L132").Value = 49788.855
$ws.Range("M132").Value = -93755798
$ws.Range("N132").Value = -54848.855

